# "Did all the excercies" - append two new glossary rows to the end of
# the (only) table, reproducing the <w:proofErr> spell-check markers
# Word leaves around words it doesn't recognise, and move the trailing
# "_GoBack" bookmark (Word's "last edit position" marker) from the old
# last cell to the new last cell.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---- Remove the "_GoBack" bookmark from its current location (end of
# the "... inerfaces implementeren" cell). It is Word's hidden
# last-edit-position bookmark, not reachable through the normal
# Bookmarks collection/Delete, so re-write that paragraph's content
# (identical text/formatting) without the bookmark tags. -------------
$goBack = $d.Content
[void]$goBack.Find.Execute("Kan van meerdere inerfaces implementeren", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$oldParaXml = $pkgOpen + `
  '<w:p w:rsidR="006240B2" w:rsidRDefault="006240B2" w:rsidP="005858BB">' + `
    '<w:r><w:t xml:space="preserve">Kan van meerdere </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>inerfaces</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> implementeren</w:t></w:r>' + `
  '</w:p>' + $pkgClose
[void]$goBack.InsertXML($oldParaXml)

# ---- Row: "CTRL +."  |  "Overriden automatisch" ----------------------
$row1 = $t.Rows.Add()
$r1idx = $t.Rows.Count

$cell1a = $t.Cell($r1idx, 1)
$cell1a.Range.Text = "CTRL +."

$cell1bXml = $pkgOpen + `
  '<w:p>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Overriden</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> automatisch</w:t></w:r>' + `
  '</w:p>' + $pkgClose
$cell1b = $t.Cell($r1idx, 2)
[void]$cell1b.Range.InsertXML($cell1bXml)

# ---- Row: "Canvas.FindName()"  |  "Gets de naam van ... for lus." ----
$row2 = $t.Rows.Add()
$r2idx = $t.Rows.Count

$cell2aXml = $pkgOpen + `
  '<w:p>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Canvas.FindName</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>()</w:t></w:r>' + `
  '</w:p>' + $pkgClose
$cell2a = $t.Cell($r2idx, 1)
[void]$cell2a.Range.InsertXML($cell2aXml)

# Last cell also carries the relocated "_GoBack" bookmark, right after
# its final run - same spot Word leaves it after the last edit.
$cell2bXml = $pkgOpen + `
  '<w:p>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Gets</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> de naam van de canvas voor bv 5 canvassen te vullen in een </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>for</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> lus.</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
  '</w:p>' + $pkgClose
$cell2b = $t.Cell($r2idx, 2)
[void]$cell2b.Range.InsertXML($cell2bXml)
